$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$elems = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates ---
# Version
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"
# Date
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"
# Contact
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet updates ---
# Definition for the root EIVL_TS element (row 2)
$elems.Range("M2").Value = "A quantity specifying a point on the axis of natural time. A point in time is most often represented as a calendar expression."
# Binding Value Set for EIVL_TS.operator (row 5)
$elems.Range("Z5").Value = "http://hl7.org/cda/stds/core/ValueSet/CDASetOperator"
# Column Z (Binding Value Set) widened to fit the new longest value (closest width reachable in this runtime's quantized column-width model)
$elems.Columns.Item(26).ColumnWidth = 50.3333333
